$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '41.530.22'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.481.48'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.27'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '92.64'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -2.45%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.86%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.510'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.98%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '32.81'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -2.45%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.862.40'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '16.41'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +10.13%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.89'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -1.96%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.412.81'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.97%  '
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.67%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '41.537.71'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.54'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0946'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '72.19'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +5.13%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '11.22'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '236.53'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '24.77'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +2.40%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.70'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '35.99'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -2.18%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '157.92'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +3.95%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.92%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0758'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '17.51'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +2.90%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.40'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -9.60%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +3.84%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -5.10%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -3.51%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.18%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.08'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -5.27%  '
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.972.38'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -1.03%  '
$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '19.25'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -3.75%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -3.08%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +1.58%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.719.33'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.05%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '68.28'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -2.54%  '
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -3.53%  '
